$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 77.5
$ws.Range("I33").Value = 77.5
$ws.Range("K33").Value = 77.5
$ws.Range("M33").Value = 151.5
$ws.Range("H43").Value = 525
$ws.Range("H53").Value = 964.86664
$ws.Range("J53").Value = 1547.7142
$ws.Range("L53").Value = 1547.7142
$ws.Range("N53").Value = -2821.7142
$ws.Range("H96").Value = 125000220
$ws.Range("I96").Value = 125000220
$ws.Range("K96").Value = 375000660
$ws.Range("M96").Value = -374999287
$ws.Range("H132").Value = 10428
$ws.Range("I132").Value = 10428
$ws.Range("K132").Value = 31284
$ws.Range("M132").Value = -28754
$ws.Range("H137").Value = 2242.652
$ws.Range("I137").Value = 1470
$ws.Range("K137").Value = 4410
$ws.Range("M137").Value = -1860
$ws.Range("H138").Value = 4896.909
$ws.Range("I138").Value = 499.5
$ws.Range("J138").Value = 7409.7144
$ws.Range("K138").Value = 1498.5
$ws.Range("L138").Value = 22229.1432
$ws.Range("M138").Value = 3641.5
$ws.Range("N138").Value = -32509.1432

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2824.5386
$ws.Range("I45").Value = 2187.5
$ws.Range("K45").Value = 2187.5
$ws.Range("M45").Value = -1810.5
$ws.Range("H61").Value = 4319.778
$ws.Range("I61").Value = 2175.6
$ws.Range("K61").Value = 2175.6
$ws.Range("M61").Value = -1963.6
$ws.Range("H74").Value = 3082.889
$ws.Range("I74").Value = 2456.375
$ws.Range("K74").Value = 2456.375
$ws.Range("M74").Value = -1582.375
$ws.Range("H77").Value = 3082.889
$ws.Range("I77").Value = 2456.375
$ws.Range("K77").Value = 12281.875
$ws.Range("M77").Value = -7913.875
$ws.Range("H110").Value = 166671250
$ws.Range("I110").Value = 333337500
$ws.Range("J110").Value = 4983.3335
$ws.Range("K110").Value = 333337500
$ws.Range("L110").Value = 4983.3335
$ws.Range("M110").Value = -333335455
$ws.Range("N110").Value = -9073.333500000001
$ws.Range("H132").Value = 1386.3334
$ws.Range("I132").Value = 1317.4375
$ws.Range("J132").Value = 1937.5
$ws.Range("K132").Value = 3952.3125
$ws.Range("L132").Value = 5812.5
$ws.Range("M132").Value = -1422.3125
$ws.Range("N132").Value = -10872.5
$ws.Range("H136").Value = 4319.778
$ws.Range("I136").Value = 2175.6
$ws.Range("K136").Value = 6526.799999999999
$ws.Range("M136").Value = -3976.799999999999

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 35719492
$ws.Range("I107").Value = 62501612
$ws.Range("K107").Value = 62501612
$ws.Range("M107").Value = -62499692
$ws.Range("H134").Value = 4653.2354
$ws.Range("I134").Value = 1133.6666
$ws.Range("J134").Value = 31050
$ws.Range("K134").Value = 3400.9998
$ws.Range("L134").Value = 93150
$ws.Range("M134").Value = -865.9998000000001
$ws.Range("N134").Value = -98220

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5277.636
$ws.Range("I31").Value = 2583.1
$ws.Range("K31").Value = 2583.1
$ws.Range("M31").Value = -2288.1
$ws.Range("H34").Value = 5277.636
$ws.Range("I34").Value = 2583.1
$ws.Range("K34").Value = 2583.1
$ws.Range("M34").Value = -2381.1
$ws.Range("H36").Value = 2048.6
$ws.Range("I36").Value = 2048.6
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 2048.6
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -1660.6
$ws.Range("N36").ClearContents()
$ws.Range("H40").Value = 2048.6
$ws.Range("I40").Value = 2048.6
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 2048.6
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -1888.6
$ws.Range("N40").ClearContents()
$ws.Range("H58").Value = 1805.0834
$ws.Range("I58").Value = 929.6111
$ws.Range("J58").Value = 4431.5
$ws.Range("K58").Value = 929.6111
$ws.Range("L58").Value = 4431.5
$ws.Range("M58").Value = -726.6111
$ws.Range("N58").Value = -4837.5
$ws.Range("H99").Value = 3750
$ws.Range("J99").Value = 5000
$ws.Range("L99").Value = 5000
$ws.Range("N99").Value = -7996
$ws.Range("H126").Value = 3750
$ws.Range("J126").Value = 5000
$ws.Range("L126").Value = 15000
$ws.Range("N126").Value = -19940
$ws.Range("H132").Value = 1150.6666
$ws.Range("I132").Value = 1150.6666
$ws.Range("K132").Value = 3451.9998
$ws.Range("M132").Value = -921.9998000000001
$ws.Range("H136").Value = 1805.0834
$ws.Range("I136").Value = 929.6111
$ws.Range("J136").Value = 4431.5
$ws.Range("K136").Value = 2788.8333
$ws.Range("L136").Value = 13294.5
$ws.Range("M136").Value = -238.8332999999998
$ws.Range("N136").Value = -18394.5

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 852
$ws.Range("I14").Value = 852
$ws.Range("K14").Value = 2556
$ws.Range("M14").Value = -2383
$ws.Range("H34").Value = 843.3333
$ws.Range("J34").Value = 944.1539
$ws.Range("L34").Value = 2832.4617
$ws.Range("N34").Value = -3000.4617

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3607.5
$ws.Range("J126").Value = 3605
$ws.Range("L126").Value = 10815
$ws.Range("N126").Value = -15755
$ws.Range("H132").Value = 54769
$ws.Range("I132").Value = 54769
$ws.Range("K132").Value = 164307
$ws.Range("M132").Value = -161777

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6125
$ws.Range("I7").Value = 5000
$ws.Range("K7").Value = 5000
$ws.Range("M7").Value = -4888
$ws.Range("H40").Value = 5386.8823
$ws.Range("I40").Value = 5155.857
$ws.Range("J40").Value = 6465
$ws.Range("K40").Value = 5155.857
$ws.Range("L40").Value = 6465
$ws.Range("M40").Value = -5019.857
$ws.Range("N40").Value = -6737
$ws.Range("H53").Value = 10000
$ws.Range("J53").Value = 10000
$ws.Range("L53").Value = 10000
$ws.Range("N53").Value = -11036
$ws.Range("H61").Value = 125005410
$ws.Range("I61").Value = 250003700
$ws.Range("K61").Value = 250003700
$ws.Range("M61").Value = -250003498
$ws.Range("H113").Value = 125005410
$ws.Range("I113").Value = 250003700
$ws.Range("K113").Value = 250003700
$ws.Range("M113").Value = -250001530
$ws.Range("H122").Value = 5000
$ws.Range("I122").Value = 5000
$ws.Range("K122").Value = 15000
$ws.Range("M122").Value = -12550
$ws.Range("H126").Value = 6125
$ws.Range("I126").Value = 5000
$ws.Range("K126").Value = 15000
$ws.Range("M126").Value = -12530
$ws.Range("H132").Value = 3882.3333
$ws.Range("I132").Value = 3879
$ws.Range("K132").Value = 11637
$ws.Range("M132").Value = -9107

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 40000
$ws.Range("J47").Value = 40000
$ws.Range("L47").Value = 40000
$ws.Range("N47").Value = -41144
$ws.Range("H51").Value = 11952.8
$ws.Range("J51").Value = 2600
$ws.Range("L51").Value = 2600
$ws.Range("N51").Value = -3620
$ws.Range("H100").Value = 1975.5
$ws.Range("I100").Value = 2334
$ws.Range("J100").Value = 900
$ws.Range("K100").Value = 4668
$ws.Range("L100").Value = 1800
$ws.Range("M100").Value = -4127
$ws.Range("N100").Value = -2882
$ws.Range("H132").Value = 677.5263
$ws.Range("I132").Value = 677.5263
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 2032.5789
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = 497.4211
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 2801.16
$ws.Range("I136").Value = 1865.5333
$ws.Range("J136").Value = 4204.6
$ws.Range("K136").Value = 5596.5999
$ws.Range("L136").Value = 12613.8
$ws.Range("M136").Value = -3046.5999
$ws.Range("N136").Value = -17713.8
